$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates derived from the crypto-price refresh diff.
# D-column price values are numeric-looking text (e.g. "42.509.91", "0.0900")
# so they are written with a leading apostrophe to force text storage and
# preserve exact formatting (matches the source inlineStr cells).
$ws.Range("D2").Value = "'42.509.91"
$ws.Range("E2").Value = "  -1.73%  "
$ws.Range("D3").Value = "'2.366.62"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").Value = "'330.81"
$ws.Range("E5").Value = "  +5.58%  "
$ws.Range("D6").Value = "'99.85"
$ws.Range("E6").Value = "  -9.04%  "
$ws.Range("E7").Value = "  -0.53%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "'0.629"
$ws.Range("E9").Value = "  -0.70%  "
$ws.Range("D10").Value = "'39.54"
$ws.Range("E10").Value = "  -7.75%  "
$ws.Range("E11").Value = "  -1.85%  "
$ws.Range("D12").Value = "'8.47"
$ws.Range("E12").Value = "  -5.34%  "
$ws.Range("E13").Value = "  -4.37%  "
$ws.Range("E14").Value = "  +0.04%  "
$ws.Range("D15").Value = "'16.34"
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("D16").Value = "'2.718.60"
$ws.Range("E16").Value = "  +0.44%  "
$ws.Range("D17").Value = "'2.353.36"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "'7.95"
$ws.Range("E18").Value = "  +8.79%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "'42.517.24"
$ws.Range("E19").Value = "  -1.75%  "
$ws.Range("E20").Value = "  -1.97%  "
$ws.Range("D21").Value = "'3.77"
$ws.Range("E21").Value = "  +9.25%  "
$ws.Range("D22").Value = "'75.93"
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("D23").Value = "'268.92"
$ws.Range("E23").Value = "  +5.76%  "
$ws.Range("E24").Value = "  -11.11%  "
$ws.Range("D25").Value = "'10.03"
$ws.Range("E25").Value = "  +9.89%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "'11.49"
$ws.Range("E27").Value = "  -4.57%  "
$ws.Range("D28").Value = "'23.22"
$ws.Range("E28").Value = "  +3.24%  "
$ws.Range("D29").Value = "'2.21"
$ws.Range("E29").Value = "  -2.29%  "
$ws.Range("D30").Value = "'176.33"
$ws.Range("E30").Value = "  +0.91%  "
$ws.Range("D31").Value = "'3.08"
$ws.Range("E31").Value = "  -2.90%  "
$ws.Range("D32").Value = "'0.0900"
$ws.Range("E32").Value = "  -3.71%  "
$ws.Range("D33").Value = "'35.35"
$ws.Range("E33").Value = "  -10.40%  "
$ws.Range("D34").Value = "'6.12"
$ws.Range("E34").Value = "  +2.20%  "
$ws.Range("E35").Value = "  -0.40%  "
$ws.Range("D36").Value = "'4.60"
$ws.Range("E36").Value = "  -7.66%  "
$ws.Range("D37").Value = "'2.98"
$ws.Range("E37").Value = "  +10.30%  "
$ws.Range("D38").Value = "'0.0359"
$ws.Range("E38").Value = "  -5.18%  "
$ws.Range("E39").Value = "  +1.28%  "
$ws.Range("E40").Value = "  -8.20%  "
$ws.Range("D41").Value = "'1.53"
$ws.Range("E41").Value = "  +3.06%  "
$ws.Range("E42").Value = "  +0.50%  "
$ws.Range("D43").Value = "'70.02"
$ws.Range("E43").Value = "  -3.85%  "
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").Value = "'118.22"
$ws.Range("E45").Value = "  +6.51%  "
$ws.Range("D46").Value = "'90.42"
$ws.Range("E46").Value = "  +28.90%  "
$ws.Range("D47").Value = "'11.89"
$ws.Range("E47").Value = "  -8.82%  "
$ws.Range("D48").Value = "'5.48"
$ws.Range("E48").Value = "  -3.10%  "
$ws.Range("D49").Value = "'9.13"
$ws.Range("E49").Value = "  -2.20%  "
$ws.Range("E50").Value = "  -2.75%  "
$ws.Range("D51").Value = "'1.568.67"
$ws.Range("E51").Value = "  +5.21%  "
